# BOT; UPDATE DATA
# Appends one more day (2020-05-19, serial 43970) of data to the three
# daily-tally sheets ("all", "kobe", "other"), pushing their trailing
# footnote row down by one, and switches the active tab from "all" to
# "kobe" (mirrors the automated daily-refresh commit).

$wb = $excel.ActiveWorkbook

$wsAll   = $wb.Worksheets.Item("all")
$wsKobe  = $wb.Worksheets.Item("kobe")
$wsOther = $wb.Worksheets.Item("other")

# ---------------------------------------------------------------------
# "all" sheet: insert new row 42 (date 43970) above the footnote row,
# which becomes row 43. Column B (累計 carried from "kobe") is left
# blank for the new day, just like the rows already on the sheet.
# ---------------------------------------------------------------------
$wsAll.Rows.Item(42).Insert() | Out-Null
$wsAll.Range("A42").Value = 43970
$wsAll.Range("C42").Value = 281
$wsAll.Range("D42").Value = 52
$wsAll.Range("E42").Value = 45
$wsAll.Range("F42").Value = 7
$wsAll.Range("G42").Value = 11
$wsAll.Range("H42").Value = 218

# ---------------------------------------------------------------------
# "kobe" sheet: insert new row 97 (date 43970) above the footnote row,
# which becomes row 98. Column B has no figure yet for the new day, so
# it is cleared outright (no cell at all, same as the source diff).
# ---------------------------------------------------------------------
$wsKobe.Rows.Item(97).Insert() | Out-Null
$wsKobe.Range("A97").Value = 43970
$wsKobe.Range("B97").Clear()
$wsKobe.Range("C97").Value = 2896
$wsKobe.Range("D97").Value = 0
$wsKobe.Range("E97").Value = 283
$wsKobe.Range("F97").Value = 47
$wsKobe.Range("G97").Value = 41
$wsKobe.Range("H97").Value = 6
$wsKobe.Range("I97").Value = 11
$wsKobe.Range("J97").Value = 209

# ---------------------------------------------------------------------
# "other" sheet: insert new row 72 (date 43970) above the footnote row,
# which becomes row 73. Figures are unchanged from the prior day.
# ---------------------------------------------------------------------
$wsOther.Rows.Item(72).Insert() | Out-Null
$wsOther.Range("A72").Value = 43970
$wsOther.Range("B72").Value = 0
$wsOther.Range("C72").Value = 14
$wsOther.Range("D72").Value = 5
$wsOther.Range("E72").Value = 4
$wsOther.Range("F72").Value = 1
$wsOther.Range("G72").Value = 0
$wsOther.Range("H72").Value = 9

# ---------------------------------------------------------------------
# View state: selections move onto the freshly entered rows, and the
# "kobe" tab becomes the active one (was "all").
# ---------------------------------------------------------------------
$wsAll.Range("I40").Select() | Out-Null
$wsOther.Range("E76").Select() | Out-Null
$wsKobe.Range("I75").Select() | Out-Null
$wsKobe.Activate() | Out-Null

Write-Host "done"
